# Add a new "GITHUB LINK" slide (Title Slide layout) as the final slide (#15).
$p = $ppt.ActivePresentation

# ppLayoutTitle = 1  -> uses the "Title Slide" custom layout (ctrTitle / subTitle placeholders),
# the same layout already used by slides 1, 5, 10, 11, 12, 13.
$s = $p.Slides.Add($p.Slides.Count + 1, 1)

# --- Title placeholder ("Title 1" / ctrTitle) ---
$title = $s.Shapes.Item(1)
$title.Name = "Title 1"
$title.TextFrame.TextRange.Text = "GITHUB LINK"

# Exact position/size (EMU 685800,533400 / 6466205x560705 expressed in points = EMU/12700)
$title.Left = 54.0
$title.Top = 42.0
$title.Width = 509.1500244140625
$title.Height = 44.150001525878906

# Disable autofit-to-shape so the body keeps a fixed size (<a:noAutofit/>)
$title.TextFrame.AutoSize = 0

# --- Subtitle placeholder ("Subtitle 2" / subTitle) ---
$subtitle = $s.Shapes.Item(2)
$subtitle.Name = "Subtitle 2"
$subtitle.TextFrame.TextRange.Text = "https://github.com/3Preeti/keylogger_newproject.git"

# Exact position/size (EMU 685800,2438400 / 8534400x276860 expressed in points)
$subtitle.Left = 54.0
$subtitle.Top = 192.0
$subtitle.Width = 672.0
$subtitle.Height = 21.80000114440918

# Hyperlink the whole subtitle run to the GitHub repo URL
$action = $subtitle.TextFrame.TextRange.ActionSettings.Item(1)
$action.Hyperlink.Address = "https://github.com/3Preeti/keylogger_newproject.git"

Write-Output $p.Slides.Count
